# Update countries & provincias Spain
#
# This script applies the "paises.xlsx" data refresh:
#  - Updates the "Datos actualizados..." timestamp string.
#  - Refreshes the daily COVID statistics (new cases, active cases, recovered,
#    critical cases, deaths) for a handful of rows following an upstream
#    source update.
#  - Corrects the ordering of a few country names that had been listed out of
#    order (Irlanda/Moldavia, Uganda/Liberia, Islas Caimanes/Trinidad y
#    Tobago, Santa Lucia/Timor Oriental) by swapping the country name text in
#    place while leaving the rest of the table layout untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Last updated timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 6 de Agosto de 2020 a las 17:27"

# --- Row 4: Estados Unidos --------------------------------------------------
$ws.Range("B4").Value = 4979568
$ws.Range("C4").Value = 6000
$ws.Range("D4").Value = 2541715
$ws.Range("E4").Value = 2276084
$ws.Range("G4").Value = 168
$ws.Range("H4").Value = 161769

# --- Row 5: Brasil -----------------------------------------------------------
$ws.Range("B5").Value = 2865053
$ws.Range("C5").Value = 2292
$ws.Range("E5").Value = 746976
$ws.Range("G5").Value = 22
$ws.Range("H5").Value = 97440

# --- Row 6: India -------------------------------------------------------------
$ws.Range("B6").Value = 1996478
$ws.Range("C6").Value = 33239
$ws.Range("D6").Value = 1351987
$ws.Range("E6").Value = 603393
$ws.Range("G6").Value = 359
$ws.Range("H6").Value = 41098

# --- Row 15: Reino Unido -------------------------------------------------------
$ws.Range("B15").Value = 308134
$ws.Range("C15").Value = 950

# --- Row 22: Alemania ----------------------------------------------------------
$ws.Range("D22").Value = 196200
$ws.Range("E22").Value = 9046

# --- Row 46: Singapur ----------------------------------------------------------
$ws.Range("D46").Value = 48031
$ws.Range("E46").Value = 6497

# --- Rows 64/65: Irlanda <-> Moldavia swap + refreshed figures ------------------
$ws.Range("A64").Value = "Moldavia"
$ws.Range("B64").Value = 26628
$ws.Range("C64").Value = 406
$ws.Range("D64").Value = 18676
$ws.Range("E64").Value = 7124
$ws.Range("G64").Value = 5
$ws.Range("H64").Value = 828

$ws.Range("A65").Value = "Irlanda"
$ws.Range("B65").Value = 26303
$ws.Range("D65").Value = 23364
$ws.Range("E65").Value = 1176
$ws.Range("H65").Value = 1763

# --- Row 73: Costa Rica ----------------------------------------------------------
$ws.Range("B73").Value = 19126
$ws.Range("C73").Value = 425
$ws.Range("D73").Value = 9236
$ws.Range("E73").Value = 9377

# --- Row 86: Cuba ------------------------------------------------------------------
$ws.Range("B86").Value = 9444
$ws.Range("C86").Value = 35
$ws.Range("E86").Value = 331

# --- Row 99: Albania ---------------------------------------------------------------
$ws.Range("B99").Value = 6016
$ws.Range("C99").Value = 127
$ws.Range("D99").Value = 3155
$ws.Range("E99").Value = 2673
$ws.Range("G99").Value = 6
$ws.Range("H99").Value = 188

# --- Row 120: Somalia --------------------------------------------------------------
$ws.Range("B120").Value = 2775
$ws.Range("C120").Value = 49
$ws.Range("D120").Value = 2409
$ws.Range("E120").Value = 278

# --- Rows 143/144: Uganda <-> Liberia swap + refreshed figures ----------------------
$ws.Range("A143").Value = "Liberia"
$ws.Range("B143").Value = 1224
$ws.Range("C143").Value = 3
$ws.Range("D143").Value = 705
$ws.Range("E143").Value = 441
$ws.Range("H143").Value = 78

$ws.Range("A144").Value = "Uganda"
$ws.Range("B144").Value = 1223
$ws.Range("C144").Value = 10
$ws.Range("D144").Value = 1102
$ws.Range("E144").Value = 116
$ws.Range("H144").Value = 5

# --- Rows 178/179: Islas Caimanes <-> Trinidad yTobago swap + refreshed figures -----
$ws.Range("A178").Value = "Trinidad yTobago"
$ws.Range("B178").Value = 207
$ws.Range("C178").Value = 8
$ws.Range("D178").Value = 135
$ws.Range("E178").Value = 64
$ws.Range("H178").Value = 8

$ws.Range("A179").Value = "Islas Caimanes"
$ws.Range("B179").Value = 203
$ws.Range("D179").Value = 202
$ws.Range("E179").Value = 0
$ws.Range("H179").Value = 1

# --- Rows 202/203: Santa Lucia <-> Timor Oriental swap (figures unchanged) ---------
$ws.Range("A202").Value = "Timor Oriental"
$ws.Range("A203").Value = "Santa Lucia"
